$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -261.2623470979252
$ws.Range("C2").Value = 25.00305362159521
$ws.Range("D2").Value = 1.000977517106549
$ws.Range("E2").Value = 8187

$ws.Range("B3").Value = -137.6682715531082
$ws.Range("C3").Value = 15.50207671634498
$ws.Range("D3").Value = 1.184165232358003
$ws.Range("E3").Value = 8186

$ws.Range("B4").Value = -61.28028343906858
$ws.Range("C4").Value = 16.67684789248625
$ws.Range("D4").Value = 1.20873786407767
$ws.Range("E4").Value = 8185

$ws.Range("B5").Value = -76.10042782040847
$ws.Range("C5").Value = 14.29618768328446
$ws.Range("D5").Value = 1.131147540983606
$ws.Range("E5").Value = 8184

$ws.Range("B6").Value = -35.24879017105882
$ws.Range("C6").Value = 12.50152755713064
$ws.Range("D6").Value = 1.096311475409836
$ws.Range("E6").Value = 8183

$ws.Range("B7").Value = 67.23491872163059
$ws.Range("C7").Value = 11.12197506722073
$ws.Range("D7").Value = 1.309644670050761
$ws.Range("E7").Value = 8182

$ws.Range("B8").Value = 122.701209946468
$ws.Range("C8").Value = 9.375381982642708
$ws.Range("D8").Value = 1.229651162790698
$ws.Range("E8").Value = 8181

$ws.Range("B9").Value = 105.472050336403
$ws.Range("C9").Value = 8.117359413202934
$ws.Range("D9").Value = 1.169934640522876
$ws.Range("E9").Value = 8180

$ws.Range("B10").Value = 31.61919792766958
$ws.Range("C10").Value = 2.579777478909402
$ws.Range("D10").Value = 1.482352941176471
$ws.Range("E10").Value = 8179

$ws.Range("B11").Value = 131.9218691300824
$ws.Range("C11").Value = 7.666911225238445
$ws.Range("D11").Value = 1.348314606741573
$ws.Range("E11").Value = 8178

$ws.Range("B12").Value = 126.5105973613434
$ws.Range("C12").Value = 7.154213036565977
$ws.Range("D12").Value = 1.25
$ws.Range("E12").Value = 8177

$ws.Range("B13").Value = 13.35291063391983
$ws.Range("C13").Value = 0.3669275929549902
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 8176

$ws.Range("B14").Value = 105.2978088406124
$ws.Range("C14").Value = 6.250764525993883
$ws.Range("D14").Value = 1.174468085106383
$ws.Range("E14").Value = 8175

$ws.Range("B15").Value = 103.214729105769
$ws.Range("C15").Value = 0.8319060435527282
$ws.Range("D15").Value = 1.833333333333333
$ws.Range("E15").Value = 8174

$ws.Range("B16").Value = 76.88872232840116
$ws.Range("C16").Value = 3.009910681512296
$ws.Range("D16").Value = 1.157894736842105
$ws.Range("E16").Value = 8173

$ws.Range("B17").Value = 93.6499293618257
$ws.Range("C17").Value = 4.221732745961821
$ws.Range("D17").Value = 1.24025974025974
$ws.Range("E17").Value = 8172

